$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.01
$ws.Range("C4").Value = -13.318
$ws.Range("D6").Value = -8.001000000000001
$ws.Range("C7").Value = -13.138
$ws.Range("D7").Value = -7.519
$ws.Range("C8").Value = -12.672
$ws.Range("D8").Value = -7.858
$ws.Range("A11").Value = -21.803
$ws.Range("A12").Value = -21.807
$ws.Range("C12").Value = -13.232
$ws.Range("C14").Value = -12.081
$ws.Range("A15").Value = -21.018
$ws.Range("D19").Value = -7.790999999999999
$ws.Range("D21").Value = -7.858
$ws.Range("C22").Value = -13.005
$ws.Range("D24").Value = -8.069999999999999
$ws.Range("D25").Value = -7.861
